$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 5201.7407
$ws.Range("I86").Value = 5431.227
$ws.Range("J86").Value = 4192
$ws.Range("K86").Value = 5431.227
$ws.Range("L86").Value = 4192
$ws.Range("M86").Value = -4308.227
$ws.Range("N86").Value = -6438

$ws.Range("H89").Value = 5201.7407
$ws.Range("I89").Value = 5431.227
$ws.Range("J89").Value = 4192
$ws.Range("K89").Value = 27156.135
$ws.Range("L89").Value = 20960
$ws.Range("M89").Value = -21540.135
$ws.Range("N89").Value = -32192

$ws.Range("H92").Value = 976.8570999999999
$ws.Range("I92").Value = 356.63635
$ws.Range("J92").Value = 3251
$ws.Range("K92").Value = 356.63635
$ws.Range("L92").Value = 3251
$ws.Range("M92").Value = 891.36365
$ws.Range("N92").Value = -5747

$ws.Range("H99").Value = 222.14285
$ws.Range("I99").Value = 256.25
$ws.Range("J99").Value = 176.66667
$ws.Range("K99").Value = 768.75
$ws.Range("L99").Value = 530.00001
$ws.Range("M99").Value = 729.25
$ws.Range("N99").Value = -3526.00001

$ws.Range("H116").Value = 4326.59
$ws.Range("I116").Value = 4960.5625
$ws.Range("J116").Value = 3885.5652
$ws.Range("K116").Value = 4960.5625
$ws.Range("L116").Value = 3885.5652
$ws.Range("M116").Value = -1518.5625
$ws.Range("N116").Value = -10769.5652

$ws.Range("H132").Value = 2675.205
$ws.Range("I132").Value = 2418.139
$ws.Range("K132").Value = 7254.417
$ws.Range("M132").Value = -4724.417

$ws.Range("H137").Value = 2055733.4
$ws.Range("I137").Value = 820555.25
$ws.Range("J137").Value = 8334556
$ws.Range("K137").Value = 2461665.75
$ws.Range("L137").Value = 25003668
$ws.Range("M137").Value = -2459115.75
$ws.Range("N137").Value = -25008768

$ws.Range("H138").Value = 1979.7833
$ws.Range("I138").Value = 1083.925
$ws.Range("J138").Value = 3771.5
$ws.Range("K138").Value = 3251.775
$ws.Range("L138").Value = 11314.5
$ws.Range("M138").Value = 1888.225
$ws.Range("N138").Value = -21594.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3822.38
$ws.Range("I32").Value = 3384.1538
$ws.Range("J32").Value = 8253.333000000001
$ws.Range("K32").Value = 3384.1538
$ws.Range("L32").Value = 8253.333000000001
$ws.Range("M32").Value = -3097.1538
$ws.Range("N32").Value = -8827.333000000001

$ws.Range("H74").Value = 728.1579
$ws.Range("I74").Value = 728.65674
$ws.Range("J74").Value = 724.44446
$ws.Range("K74").Value = 728.65674
$ws.Range("L74").Value = 724.44446
$ws.Range("M74").Value = 145.34326
$ws.Range("N74").Value = -2472.44446

$ws.Range("H77").Value = 728.1579
$ws.Range("I77").Value = 728.65674
$ws.Range("J77").Value = 724.44446
$ws.Range("K77").Value = 3643.2837
$ws.Range("L77").Value = 3622.2223
$ws.Range("M77").Value = 724.7163
$ws.Range("N77").Value = -12358.2223

$ws.Range("H124").Value = 20666.479
$ws.Range("J124").Value = 20666.479
$ws.Range("L124").Value = 20666.479
$ws.Range("N124").Value = -30486.479

$ws.Range("H132").Value = 64978.09
$ws.Range("I132").Value = 81817.664
$ws.Range("J132").Value = 3563.1765
$ws.Range("K132").Value = 245452.992
$ws.Range("L132").Value = 10689.5295
$ws.Range("M132").Value = -242922.992
$ws.Range("N132").Value = -15749.5295

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1653.7084
$ws.Range("I31").Value = 1364.9429
$ws.Range("J31").Value = 2431.1538
$ws.Range("K31").Value = 1364.9429
$ws.Range("L31").Value = 2431.1538
$ws.Range("M31").Value = -1069.9429
$ws.Range("N31").Value = -3021.1538

$ws.Range("H34").Value = 1653.7084
$ws.Range("I34").Value = 1364.9429
$ws.Range("J34").Value = 2431.1538
$ws.Range("K34").Value = 1364.9429
$ws.Range("L34").Value = 2431.1538
$ws.Range("M34").Value = -1162.9429
$ws.Range("N34").Value = -2835.1538

$ws.Range("H58").Value = 969.8182
$ws.Range("I58").Value = 943.95
$ws.Range("J58").Value = 1228.5
$ws.Range("K58").Value = 943.95
$ws.Range("L58").Value = 1228.5
$ws.Range("M58").Value = -740.95
$ws.Range("N58").Value = -1634.5

$ws.Range("H99").Value = 2061.111
$ws.Range("I99").Value = 2068.75
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 2068.75
$ws.Range("L99").Value = 2000
$ws.Range("M99").Value = -570.75
$ws.Range("N99").Value = -4996

$ws.Range("H126").Value = 2061.111
$ws.Range("I126").Value = 2068.75
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 6206.25
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -3736.25
$ws.Range("N126").Value = -10940

$ws.Range("H132").Value = 3444.6667
$ws.Range("I132").Value = 2666.8845
$ws.Range("J132").Value = 8500.25
$ws.Range("K132").Value = 8000.6535
$ws.Range("L132").Value = 25500.75
$ws.Range("M132").Value = -5470.6535
$ws.Range("N132").Value = -30560.75

$ws.Range("H134").Value = 2767.3572
$ws.Range("I134").Value = 2859.4211
$ws.Range("J134").Value = 2363.6924
$ws.Range("K134").Value = 8578.263300000001
$ws.Range("L134").Value = 7091.0772
$ws.Range("M134").Value = -6043.263300000001
$ws.Range("N134").Value = -12161.0772

$ws.Range("H136").Value = 969.8182
$ws.Range("I136").Value = 943.95
$ws.Range("J136").Value = 1228.5
$ws.Range("K136").Value = 2831.85
$ws.Range("L136").Value = 3685.5
$ws.Range("M136").Value = -281.8500000000004
$ws.Range("N136").Value = -8785.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 581065.2
$ws.Range("I121").Value = 600
$ws.Range("J121").Value = 636347.5600000001
$ws.Range("K121").Value = 1800
$ws.Range("L121").Value = 1909042.68
$ws.Range("M121").Value = -490
$ws.Range("N121").Value = -1911662.68

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1342.5217
$ws.Range("I113").Value = 1184.8667
$ws.Range("J113").Value = 1638.125
$ws.Range("K113").Value = 1184.8667
$ws.Range("L113").Value = 1638.125
$ws.Range("M113").Value = 985.1333
$ws.Range("N113").Value = -5978.125

$ws.Range("H126").Value = 13181.261
$ws.Range("I126").Value = 2999
$ws.Range("J126").Value = 22515
$ws.Range("K126").Value = 8997
$ws.Range("L126").Value = 67545
$ws.Range("M126").Value = -6527
$ws.Range("N126").Value = -72485

$ws.Range("H132").Value = 1573.2391
$ws.Range("I132").Value = 1188.303
$ws.Range("J132").Value = 2550.3845
$ws.Range("K132").Value = 3564.909000000001
$ws.Range("L132").Value = 7651.1535
$ws.Range("M132").Value = -1034.909000000001
$ws.Range("N132").Value = -12711.1535

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6668648
$ws.Range("I7").Value = 10001651
$ws.Range("K7").Value = 10001651
$ws.Range("M7").Value = -10001539

$ws.Range("H126").Value = 6668648
$ws.Range("I126").Value = 10001651
$ws.Range("K126").Value = 30004953
$ws.Range("M126").Value = -30002483

$ws.Range("H132").Value = 2455.575
$ws.Range("I132").Value = 2201.0557
$ws.Range("J132").Value = 4746.25
$ws.Range("K132").Value = 6603.1671
$ws.Range("L132").Value = 14238.75
$ws.Range("M132").Value = -4073.1671
$ws.Range("N132").Value = -19298.75

$ws.Range("H136").Value = 1459.1786
$ws.Range("I136").Value = 1359.2174
$ws.Range("K136").Value = 4077.6522
$ws.Range("M136").Value = -1527.6522

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 9412
$ws.Range("J122").Value = 3416.6667
$ws.Range("L122").Value = 10250.0001
$ws.Range("N122").Value = -15150.0001

$ws.Range("H132").Value = 1724.7966
$ws.Range("I132").Value = 2279.1353
$ws.Range("J132").Value = 792.5
$ws.Range("K132").Value = 6837.4059
$ws.Range("L132").Value = 2377.5
$ws.Range("M132").Value = -4307.4059
$ws.Range("N132").Value = -7437.5

$ws.Range("H136").Value = 1218.2084
$ws.Range("I136").Value = 1084.9517
$ws.Range("J136").Value = 2044.4
$ws.Range("K136").Value = 3254.8551
$ws.Range("L136").Value = 6133.200000000001
$ws.Range("M136").Value = -704.8551000000002
$ws.Range("N136").Value = -11233.2
